# Vorlage_Rechnung.docx - "Adjusted Template with new title"
#
# The small address table in the header area (table #2) has its first
# two rows reworked:
#   - Row 1 used to hold "BärFoods · Eigerstrasse 74 3007 Bern" in 8pt
#     (sz=16) text spread across many runs (with a stray _GoBack
#     bookmark in the middle). It becomes just a bold 20pt (sz=40)
#     "BärFoods" title, and the row is taller to fit it.
#   - Row 2's previously-empty paragraph now carries the _GoBack
#     bookmark that used to sit inside row 1.

$d = $word.ActiveDocument
$t = $d.Tables.Item(2)

# 1) Row 1: trHeight 265 -> 450 (twentieths of a point == 13.25pt -> 22.5pt)
$row1 = $t.Rows.Item(1)
$row1.Height = 22.5

# 2) Merge the "BärF" + bookmark + "oods" runs into a single "BärFoods"
#    run (this also drops the _GoBack bookmark that used to live between
#    them, matching the target XML).
$found1 = $d.Content.Find.Execute("BärFoods", $true, $false, $false, $false, $false, $true, 1, $false, "BärFoods", 2)

# 3) Drop the trailing address text, leaving "BärFoods" followed
#    immediately by the closing proofErr spellEnd mark.
$found2 = $d.Content.Find.Execute(" · Eigerstrasse 74 3007 Bern", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 4) Make the title bold (paragraph mark + run) and bump the run's font
#    size to 20pt (sz=40 half-points).
$cell1 = $row1.Cells.Item(1)
$para1 = $cell1.Range.Paragraphs.Item(1)
$pRange = $para1.Range
$pRange.Font.Bold = 1

$runRange = $d.Range($pRange.Start, $pRange.End - 1)
$runRange.Font.Size = 20

# 5) Move the _GoBack bookmark into row 2's (still empty) paragraph.
$row2 = $t.Rows.Item(2)
$cell2 = $row2.Cells.Item(1)
$para2 = $cell2.Range.Paragraphs.Item(1)
$d.Bookmarks.Add("_GoBack", $para2.Range) | Out-Null
